$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A few Price cells are numeric-looking strings with significant trailing
# zeros (e.g. "6.90", "1.00", "2.70"). Plain .Value assignment lets Excel
# auto-detect them as numbers and silently drop the trailing zero, so force
# a Text number format on exactly those cells first.
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"

$ws.Range("D2").Value = "64.138.74"
$ws.Range("E2").Value = "  -3.32%  "
$ws.Range("D3").Value = "3.128.11"
$ws.Range("E3").Value = "  -2.74%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "608.65"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").Value = "147.28"
$ws.Range("E6").Value = "  -5.06%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.126.73"
$ws.Range("E8").Value = "  -2.76%  "
$ws.Range("D9").Value = "0.525"
$ws.Range("E9").Value = "  -3.60%  "
$ws.Range("D10").Value = "0.151"
$ws.Range("E10").Value = "  -5.56%  "
$ws.Range("D11").Value = "5.54"
$ws.Range("E11").Value = "  -2.78%  "
$ws.Range("D12").Value = "0.474"
$ws.Range("E12").Value = "  -5.27%  "
$ws.Range("D13").Value = "0.0000255"
$ws.Range("E13").Value = "  -4.46%  "
$ws.Range("D14").Value = "36.28"
$ws.Range("E14").Value = "  -5.14%  "
$ws.Range("D15").Value = "3.642.43"
$ws.Range("E15").Value = "  -2.69%  "
$ws.Range("D16").Value = "64.099.20"
$ws.Range("E16").Value = "  -3.47%  "
$ws.Range("D18").Value = "3.136.04"
$ws.Range("E18").Value = "  -2.60%  "
$ws.Range("D19").Value = "6.90"
$ws.Range("E19").Value = "  -4.69%  "
$ws.Range("D20").Value = "477.83"
$ws.Range("E20").Value = "  -5.54%  "
$ws.Range("D21").Value = "14.49"
$ws.Range("E21").Value = "  -4.69%  "
$ws.Range("D22").Value = "0.703"
$ws.Range("E22").Value = "  -3.37%  "
$ws.Range("D23").Value = "7.68"
$ws.Range("E23").Value = "  -3.67%  "
$ws.Range("D24").Value = "13.68"
$ws.Range("E24").Value = "  -5.61%  "
$ws.Range("D25").Value = "82.97"
$ws.Range("E25").Value = "  -2.47%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  -2.74%  "
$ws.Range("D28").Value = "8.41"
$ws.Range("E28").Value = "  -6.47%  "
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  -5.91%  "
$ws.Range("D30").Value = "0.122"
$ws.Range("E30").Value = "  -18.97%  "
$ws.Range("E31").Value = "  -2.41%  "
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").Value = "2.70"
$ws.Range("E33").Value = "  -6.44%  "
$ws.Range("D34").Value = "26.29"
$ws.Range("E34").Value = "  -6.90%  "
$ws.Range("E35").Value = "  -5.78%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "6.02"
$ws.Range("E36").Value = "  -5.63%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "54.35"
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "3.06"
$ws.Range("E38").Value = "  +1.56%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0726"
$ws.Range("E39").Value = "  -5.74%  "
$ws.Range("D40").Value = "451.27"
$ws.Range("E40").Value = "  -9.67%  "
$ws.Range("D41").Value = "0.0396"
$ws.Range("E41").Value = "  -5.16%  "
$ws.Range("D42").Value = "0.123"
$ws.Range("E42").Value = "  -4.96%  "
$ws.Range("D43").Value = "8.38"
$ws.Range("E43").Value = "  -3.74%  "
$ws.Range("D44").Value = "2.855.95"
$ws.Range("E44").Value = "  -2.20%  "
$ws.Range("D45").Value = "0.269"
$ws.Range("E45").Value = "  -8.33%  "
$ws.Range("E46").Value = "  -7.96%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "26.35"
$ws.Range("E47").Value = "  -5.77%  "
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value = "0.998"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "2.31"
$ws.Range("E49").Value = "  -3.22%  "
$ws.Range("E50").Value = "  -2.63%  "
$ws.Range("D51").Value = "118.65"
$ws.Range("E51").Value = "  -2.18%  "
